# Slide 6: "마케팅/기술영업 전문 인력 모집" job-posting slide.
# The author deleted the word "전문 " from the heading textbox, which
# (via spAutoFit + manual nudging) also shifted a couple of nearby
# pictures/textboxes. Reproduce both effects explicitly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$grp = $s.Shapes.Item(1)

# EMU -> point helper (1 pt = 12700 EMU). Add a half-EMU epsilon so the
# float->EMU rounding on save lands on the exact target integer instead
# of truncating one EMU short.
function EmuToPt([double]$emu) {
    return ($emu + 0.5) / 12700
}

# --- Text edit: drop "전문 " from the "마케팅/기술영업 전문 인력 모집" heading ---
$heading = $grp.GroupItems.Item(3)   # "TextBox 4"
$found = $heading.TextFrame.TextRange.Find("전문 ")
if ($found -ne $null) {
    $found.Text = ""
}

# --- Reposition/resize shapes to match the final layout ---

# "그림 2" picture (logo strip) shifts right; size unchanged.
$pic2 = $grp.GroupItems.Item(2)
$pic2.Left = EmuToPt 2942155

# Heading textbox narrows (autofit) and shifts right; top/height unchanged.
$heading.Left = EmuToPt 2421080
$heading.Width = EmuToPt 4407686

# "그림 5" picture (small badge) moves left and slightly up.
$pic5 = $grp.GroupItems.Item(4)
$pic5.Left = EmuToPt 6776811
$pic5.Top = EmuToPt 2407212

# "TextBox 10" (date range) nudges left and up slightly.
$tb10 = $grp.GroupItems.Item(6)
$tb10.Left = EmuToPt 3739647
$tb10.Top = EmuToPt 2846267
